# Applies the "Additions to final manuscript" edit to the key-coding sheet.
#
# Net effect described by the diff:
#   - Row 13 ("has_cv2" / duplicate CNS description) is removed entirely;
#     rows 14-27 shift up by one (so former row 14 "has_cns" becomes row 13, etc.)
#   - The row that ends up at position 24 is repurposed from "has_dem_and_cva_or_degen"
#     (which now lives at row 23) to a brand-new entry "FinalTx_coll" with updated coding.
#   - The row that ends up at position 25 keeps the label "PercOSA" but gets new,
#     simplified coding text (rather than the old "FinalTx" content sliding into it).
#   - The final row (26) is "StudyType" with its original text, and the sheet now
#     only spans down to row 26 (was 27).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove row 13 ("has_cv2"); everything below shifts up by one row automatically.
$ws.Rows("13").Delete()

# After the shift, row 24 held "has_dem_and_cva_or_degen" (moved up from the old row 25's
# slot indirectly) -- overwrite it with the new FinalTx_coll entry.
$ws.Range("A24").Value = "FinalTx_coll"
$ws.Range("B24").Value = "0 = BPAP,  1 = ASV,  2 = CPAP,  3 = Other,  "

# Row 25 keeps its "PercOSA" label but the coding description is replaced.
$ws.Range("A25").Value = "PercOSA"
$ws.Range("B25").Value = "0 = mostly_OSA,  1 = mostly_CSA,  "

# Row 26 remains "StudyType" with unchanged content (already correct after the row delete).
